$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 3.554827
$ws.Cells.Item(2, 8).Value = 10.664481
$ws.Cells.Item(2, 9).Value = 0.2148969460055877
$ws.Cells.Item(2, 10).Value = 0.2148969460055877
$ws.Cells.Item(2, 13).Value = 34.293805
$ws.Cells.Item(2, 14).Value = 102.881415
$ws.Cells.Item(2, 15).Value = 0.3000373067112135
$ws.Cells.Item(2, 16).Value = 0.3000373067112135
$ws.Cells.Item(2, 17).Value = 121.908543946735
$ws.Cells.Item(2, 18).Value = 1097.176895520615
$ws.Cells.Item(2, 19).Value = 0.0644771008999816
$ws.Cells.Item(2, 20).Value = 0.0644771008999816

$ws.Cells.Item(3, 7).Value = 3.554827
$ws.Cells.Item(3, 8).Value = 10.664481
$ws.Cells.Item(3, 9).Value = 0.2148969460055877
$ws.Cells.Item(3, 10).Value = 0.2148969460055877
$ws.Cells.Item(3, 15).Value = 0.2416702056223797
$ws.Cells.Item(3, 16).Value = 0.2416702056223798
$ws.Cells.Item(3, 17).Value = 98.19333204150266
$ws.Cells.Item(3, 18).Value = 883.7399883735238
$ws.Cells.Item(3, 19).Value = 0.05193418912879182
$ws.Cells.Item(3, 20).Value = 0.05193418912879182

$ws.Cells.Item(4, 7).Value = 3.554827
$ws.Cells.Item(4, 8).Value = 10.664481
$ws.Cells.Item(4, 9).Value = 0.2148969460055877
$ws.Cells.Item(4, 10).Value = 0.2148969460055877
$ws.Cells.Item(4, 13).Value = 32.81168366666667
$ws.Cells.Item(4, 14).Value = 98.435051
$ws.Cells.Item(4, 15).Value = 0.2870701922987834
$ws.Cells.Item(4, 16).Value = 0.2870701922987835
$ws.Cells.Item(4, 17).Value = 116.6398590137257
$ws.Cells.Item(4, 18).Value = 1049.758731123531
$ws.Cells.Item(4, 19).Value = 0.06169050761424535
$ws.Cells.Item(4, 20).Value = 0.06169050761424535

$ws.Cells.Item(5, 7).Value = 3.554827
$ws.Cells.Item(5, 8).Value = 10.664481
$ws.Cells.Item(5, 9).Value = 0.2148969460055877
$ws.Cells.Item(5, 10).Value = 0.2148969460055877
$ws.Cells.Item(5, 13).Value = 6.520685
$ws.Cells.Item(5, 14).Value = 19.562055
$ws.Cells.Item(5, 15).Value = 0.05704962646496092
$ws.Cells.Item(5, 16).Value = 0.05704962646496093
$ws.Cells.Item(5, 17).Value = 23.179907096495
$ws.Cells.Item(5, 18).Value = 208.619163868455
$ws.Cells.Item(5, 19).Value = 0.01225979049807966
$ws.Cells.Item(5, 20).Value = 0.01225979049807966

$ws.Cells.Item(6, 7).Value = 3.554827
$ws.Cells.Item(6, 8).Value = 10.664481
$ws.Cells.Item(6, 9).Value = 0.2148969460055877
$ws.Cells.Item(6, 10).Value = 0.2148969460055877
$ws.Cells.Item(6, 13).Value = 13.04976133333333
$ws.Cells.Item(6, 14).Value = 39.14928399999999
$ws.Cells.Item(6, 15).Value = 0.1141726689026624
$ws.Cells.Item(6, 16).Value = 0.1141726689026624
$ws.Cells.Item(6, 17).Value = 46.38964393128932
$ws.Cells.Item(6, 18).Value = 417.5067953816039
$ws.Cells.Item(6, 19).Value = 0.02453535786448927
$ws.Cells.Item(6, 20).Value = 0.02453535786448928

$ws.Cells.Item(7, 9).Value = 0.3107709374420163
$ws.Cells.Item(7, 10).Value = 0.3107709374420163
$ws.Cells.Item(7, 13).Value = 34.293805
$ws.Cells.Item(7, 14).Value = 102.881415
$ws.Cells.Item(7, 15).Value = 0.3000373067112135
$ws.Cells.Item(7, 16).Value = 0.3000373067112135
$ws.Cells.Item(7, 17).Value = 176.2967468301433
$ws.Cells.Item(7, 18).Value = 1586.67072147129
$ws.Cells.Item(7, 19).Value = 0.09324287507422158
$ws.Cells.Item(7, 20).Value = 0.09324287507422158

$ws.Cells.Item(8, 9).Value = 0.3107709374420163
$ws.Cells.Item(8, 10).Value = 0.3107709374420163
$ws.Cells.Item(8, 15).Value = 0.2416702056223797
$ws.Cells.Item(8, 16).Value = 0.2416702056223798
$ws.Cells.Item(8, 19).Value = 0.07510407635307179
$ws.Cells.Item(8, 20).Value = 0.07510407635307179

$ws.Cells.Item(9, 9).Value = 0.3107709374420163
$ws.Cells.Item(9, 10).Value = 0.3107709374420163
$ws.Cells.Item(9, 13).Value = 32.81168366666667
$ws.Cells.Item(9, 14).Value = 98.435051
$ws.Cells.Item(9, 15).Value = 0.2870701922987834
$ws.Cells.Item(9, 16).Value = 0.2870701922987835
$ws.Cells.Item(9, 17).Value = 168.6774940387362
$ws.Cells.Item(9, 18).Value = 1518.097446348626
$ws.Cells.Item(9, 19).Value = 0.08921307277235281
$ws.Cells.Item(9, 20).Value = 0.08921307277235282

$ws.Cells.Item(10, 9).Value = 0.3107709374420163
$ws.Cells.Item(10, 10).Value = 0.3107709374420163
$ws.Cells.Item(10, 13).Value = 6.520685
$ws.Cells.Item(10, 14).Value = 19.562055
$ws.Cells.Item(10, 15).Value = 0.05704962646496092
$ws.Cells.Item(10, 16).Value = 0.05704962646496093
$ws.Cells.Item(10, 17).Value = 33.52137660443667
$ws.Cells.Item(10, 18).Value = 301.69238943993
$ws.Cells.Item(10, 19).Value = 0.01772936589723277
$ws.Cells.Item(10, 20).Value = 0.01772936589723277

$ws.Cells.Item(11, 9).Value = 0.3107709374420163
$ws.Cells.Item(11, 10).Value = 0.3107709374420163
$ws.Cells.Item(11, 13).Value = 13.04976133333333
$ws.Cells.Item(11, 14).Value = 39.14928399999999
$ws.Cells.Item(11, 15).Value = 0.1141726689026624
$ws.Cells.Item(11, 16).Value = 0.1141726689026624
$ws.Cells.Item(11, 17).Value = 67.08589116828709
$ws.Cells.Item(11, 18).Value = 603.7730205145839
$ws.Cells.Item(11, 19).Value = 0.03548154734513732
$ws.Cells.Item(11, 20).Value = 0.03548154734513733

$ws.Cells.Item(12, 7).Value = 2.615693
$ws.Cells.Item(12, 8).Value = 7.847079000000001
$ws.Cells.Item(12, 9).Value = 0.1581242736673807
$ws.Cells.Item(12, 10).Value = 0.1581242736673807
$ws.Cells.Item(12, 13).Value = 34.293805
$ws.Cells.Item(12, 14).Value = 102.881415
$ws.Cells.Item(12, 15).Value = 0.3000373067112135
$ws.Cells.Item(12, 16).Value = 0.3000373067112135
$ws.Cells.Item(12, 17).Value = 89.70206568186501
$ws.Cells.Item(12, 18).Value = 807.318591136785
$ws.Cells.Item(12, 19).Value = 0.04744318119682775
$ws.Cells.Item(12, 20).Value = 0.04744318119682775

$ws.Cells.Item(13, 7).Value = 2.615693
$ws.Cells.Item(13, 8).Value = 7.847079000000001
$ws.Cells.Item(13, 9).Value = 0.1581242736673807
$ws.Cells.Item(13, 10).Value = 0.1581242736673807
$ws.Cells.Item(13, 15).Value = 0.2416702056223797
$ws.Cells.Item(13, 16).Value = 0.2416702056223798
$ws.Cells.Item(13, 17).Value = 72.25207056985734
$ws.Cells.Item(13, 18).Value = 650.268635128716
$ws.Cells.Item(13, 19).Value = 0.03821392573108533
$ws.Cells.Item(13, 20).Value = 0.03821392573108533

$ws.Cells.Item(14, 7).Value = 2.615693
$ws.Cells.Item(14, 8).Value = 7.847079000000001
$ws.Cells.Item(14, 9).Value = 0.1581242736673807
$ws.Cells.Item(14, 10).Value = 0.1581242736673807
$ws.Cells.Item(14, 13).Value = 32.81168366666667
$ws.Cells.Item(14, 14).Value = 98.435051
$ws.Cells.Item(14, 15).Value = 0.2870701922987834
$ws.Cells.Item(14, 16).Value = 0.2870701922987835
$ws.Cells.Item(14, 17).Value = 85.82529128511435
$ws.Cells.Item(14, 18).Value = 772.4276215660291
$ws.Cells.Item(14, 19).Value = 0.04539276564880042
$ws.Cells.Item(14, 20).Value = 0.04539276564880043

$ws.Cells.Item(15, 7).Value = 2.615693
$ws.Cells.Item(15, 8).Value = 7.847079000000001
$ws.Cells.Item(15, 9).Value = 0.1581242736673807
$ws.Cells.Item(15, 10).Value = 0.1581242736673807
$ws.Cells.Item(15, 13).Value = 6.520685
$ws.Cells.Item(15, 14).Value = 19.562055
$ws.Cells.Item(15, 15).Value = 0.05704962646496092
$ws.Cells.Item(15, 16).Value = 0.05704962646496093
$ws.Cells.Item(15, 17).Value = 17.056110109705
$ws.Cells.Item(15, 18).Value = 153.504990987345
$ws.Cells.Item(15, 19).Value = 0.009020930747767324
$ws.Cells.Item(15, 20).Value = 0.009020930747767324

$ws.Cells.Item(16, 7).Value = 2.615693
$ws.Cells.Item(16, 8).Value = 7.847079000000001
$ws.Cells.Item(16, 9).Value = 0.1581242736673807
$ws.Cells.Item(16, 10).Value = 0.1581242736673807
$ws.Cells.Item(16, 13).Value = 13.04976133333333
$ws.Cells.Item(16, 14).Value = 39.14928399999999
$ws.Cells.Item(16, 15).Value = 0.1141726689026624
$ws.Cells.Item(16, 16).Value = 0.1141726689026624
$ws.Cells.Item(16, 17).Value = 34.13416937127067
$ws.Cells.Item(16, 18).Value = 307.207524341436
$ws.Cells.Item(16, 19).Value = 0.01805347034289983
$ws.Cells.Item(16, 20).Value = 0.01805347034289983

$ws.Cells.Item(17, 7).Value = 4.248598333333334
$ws.Cells.Item(17, 8).Value = 12.745795
$ws.Cells.Item(17, 9).Value = 0.2568369168563656
$ws.Cells.Item(17, 10).Value = 0.2568369168563656
$ws.Cells.Item(17, 13).Value = 34.293805
$ws.Cells.Item(17, 14).Value = 102.881415
$ws.Cells.Item(17, 15).Value = 0.3000373067112135
$ws.Cells.Item(17, 16).Value = 0.3000373067112135
$ws.Cells.Item(17, 17).Value = 145.7006027666583
$ws.Cells.Item(17, 18).Value = 1311.305424899925
$ws.Cells.Item(17, 19).Value = 0.07706065679759579
$ws.Cells.Item(17, 20).Value = 0.07706065679759579

$ws.Cells.Item(18, 7).Value = 4.248598333333334
$ws.Cells.Item(18, 8).Value = 12.745795
$ws.Cells.Item(18, 9).Value = 0.2568369168563656
$ws.Cells.Item(18, 10).Value = 0.2568369168563656
$ws.Cells.Item(18, 15).Value = 0.2416702056223797
$ws.Cells.Item(18, 16).Value = 0.2416702056223798
$ws.Cells.Item(18, 17).Value = 117.3570547472422
$ws.Cells.Item(18, 18).Value = 1056.21349272518
$ws.Cells.Item(18, 19).Value = 0.06206983050809591
$ws.Cells.Item(18, 20).Value = 0.06206983050809592

$ws.Cells.Item(19, 7).Value = 4.248598333333334
$ws.Cells.Item(19, 8).Value = 12.745795
$ws.Cells.Item(19, 9).Value = 0.2568369168563656
$ws.Cells.Item(19, 10).Value = 0.2568369168563656
$ws.Cells.Item(19, 13).Value = 32.81168366666667
$ws.Cells.Item(19, 14).Value = 98.435051
$ws.Cells.Item(19, 15).Value = 0.2870701922987834
$ws.Cells.Item(19, 16).Value = 0.2870701922987835
$ws.Cells.Item(19, 17).Value = 139.4036645400606
$ws.Cells.Item(19, 18).Value = 1254.632980860545
$ws.Cells.Item(19, 19).Value = 0.07373022311138351
$ws.Cells.Item(19, 20).Value = 0.07373022311138352

$ws.Cells.Item(20, 7).Value = 4.248598333333334
$ws.Cells.Item(20, 8).Value = 12.745795
$ws.Cells.Item(20, 9).Value = 0.2568369168563656
$ws.Cells.Item(20, 10).Value = 0.2568369168563656
$ws.Cells.Item(20, 13).Value = 6.520685
$ws.Cells.Item(20, 14).Value = 19.562055
$ws.Cells.Item(20, 15).Value = 0.05704962646496092
$ws.Cells.Item(20, 16).Value = 0.05704962646496093
$ws.Cells.Item(20, 17).Value = 27.70377142319167
$ws.Cells.Item(20, 18).Value = 249.333942808725
$ws.Cells.Item(20, 19).Value = 0.01465245016906788
$ws.Cells.Item(20, 20).Value = 0.01465245016906788

$ws.Cells.Item(21, 7).Value = 4.248598333333334
$ws.Cells.Item(21, 8).Value = 12.745795
$ws.Cells.Item(21, 9).Value = 0.2568369168563656
$ws.Cells.Item(21, 10).Value = 0.2568369168563656
$ws.Cells.Item(21, 13).Value = 13.04976133333333
$ws.Cells.Item(21, 14).Value = 39.14928399999999
$ws.Cells.Item(21, 15).Value = 0.1141726689026624
$ws.Cells.Item(21, 16).Value = 0.1141726689026624
$ws.Cells.Item(21, 17).Value = 55.44319425119777
$ws.Cells.Item(21, 18).Value = 498.98874826078
$ws.Cells.Item(21, 19).Value = 0.02932375627022244
$ws.Cells.Item(21, 20).Value = 0.02932375627022245

$ws.Cells.Item(22, 7).Value = 0.9821143333333332
$ws.Cells.Item(22, 8).Value = 2.946343
$ws.Cells.Item(22, 9).Value = 0.0593709260286498
$ws.Cells.Item(22, 10).Value = 0.0593709260286498
$ws.Cells.Item(22, 13).Value = 34.293805
$ws.Cells.Item(22, 14).Value = 102.881415
$ws.Cells.Item(22, 15).Value = 0.3000373067112135
$ws.Cells.Item(22, 16).Value = 0.3000373067112135
$ws.Cells.Item(22, 17).Value = 33.68043743503833
$ws.Cells.Item(22, 18).Value = 303.1239369153449
$ws.Cells.Item(22, 19).Value = 0.01781349274258677
$ws.Cells.Item(22, 20).Value = 0.01781349274258677

$ws.Cells.Item(23, 7).Value = 0.9821143333333332
$ws.Cells.Item(23, 8).Value = 2.946343
$ws.Cells.Item(23, 9).Value = 0.0593709260286498
$ws.Cells.Item(23, 10).Value = 0.0593709260286498
$ws.Cells.Item(23, 15).Value = 0.2416702056223797
$ws.Cells.Item(23, 16).Value = 0.2416702056223798
$ws.Cells.Item(23, 17).Value = 27.12848721913022
$ws.Cells.Item(23, 18).Value = 244.156384972172
$ws.Cells.Item(23, 19).Value = 0.01434818390133489
$ws.Cells.Item(23, 20).Value = 0.0143481839013349

$ws.Cells.Item(24, 7).Value = 0.9821143333333332
$ws.Cells.Item(24, 8).Value = 2.946343
$ws.Cells.Item(24, 9).Value = 0.0593709260286498
$ws.Cells.Item(24, 10).Value = 0.0593709260286498
$ws.Cells.Item(24, 13).Value = 32.81168366666667
$ws.Cells.Item(24, 14).Value = 98.435051
$ws.Cells.Item(24, 15).Value = 0.2870701922987834
$ws.Cells.Item(24, 16).Value = 0.2870701922987835
$ws.Cells.Item(24, 17).Value = 32.22482482983255
$ws.Cells.Item(24, 18).Value = 290.023423468493
$ws.Cells.Item(24, 19).Value = 0.01704362315200135
$ws.Cells.Item(24, 20).Value = 0.01704362315200135

$ws.Cells.Item(25, 7).Value = 0.9821143333333332
$ws.Cells.Item(25, 8).Value = 2.946343
$ws.Cells.Item(25, 9).Value = 0.0593709260286498
$ws.Cells.Item(25, 10).Value = 0.0593709260286498
$ws.Cells.Item(25, 13).Value = 6.520685
$ws.Cells.Item(25, 14).Value = 19.562055
$ws.Cells.Item(25, 15).Value = 0.05704962646496092
$ws.Cells.Item(25, 16).Value = 0.05704962646496093
$ws.Cells.Item(25, 17).Value = 6.404058201651666
$ws.Cells.Item(25, 18).Value = 57.63652381486499
$ws.Cells.Item(25, 19).Value = 0.003387089152813297
$ws.Cells.Item(25, 20).Value = 0.003387089152813297

$ws.Cells.Item(26, 7).Value = 0.9821143333333332
$ws.Cells.Item(26, 8).Value = 2.946343
$ws.Cells.Item(26, 9).Value = 0.0593709260286498
$ws.Cells.Item(26, 10).Value = 0.0593709260286498
$ws.Cells.Item(26, 13).Value = 13.04976133333333
$ws.Cells.Item(26, 14).Value = 39.14928399999999
$ws.Cells.Item(26, 15).Value = 0.1141726689026624
$ws.Cells.Item(26, 16).Value = 0.1141726689026624
$ws.Cells.Item(26, 17).Value = 12.81635765204577
$ws.Cells.Item(26, 18).Value = 115.347218868412
$ws.Cells.Item(26, 19).Value = 0.006778537079913492
$ws.Cells.Item(26, 20).Value = 0.006778537079913494
